# Natmi following Dr Hou advice
# Update LR-pair data: add FAPs cluster rows and recompute values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-10: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# then 16 numeric columns (E..T)
$data = @(
    @("ECs",  "L1cam", "Erbb3", "ECs",  3, 1, 19.72083766666667, 59.162513,         0.8016210077351786,  0.8016210077351787,  2, 0.6666666666666666, 0.2167736666666666, 0.6503209999999999, 0.06003183809373913, 0.06003183809373914, 4.274958290741444,  38.474624616673,     0.04812278254889825, 0.04812278254889826),
    @("ECs",  "L1cam", "Erbb3", "FAPs", 3, 1, 19.72083766666667, 59.162513,         0.8016210077351786,  0.8016210077351787,  2, 0.6666666666666666, 0.3460976666666666, 1.038293,            0.09584595495126666, 0.09584595495126667, 6.825335901145444,  61.428023110309,     0.07683213099537492, 0.07683213099537493),
    @("ECs",  "L1cam", "Erbb3", "sCs",  3, 1, 19.72083766666667, 59.162513,         0.8016210077351786,  0.8016210077351787,  3, 1,                   3.048107,            9.144321,            0.8441222069549942,  0.8441222069549942,  60.11122333763034,  541.001010038673,    0.6766660941909054,  0.6766660941909055),
    @("FAPs", "L1cam", "Erbb3", "ECs",  1, 0.3333333333333333, 0.099159, 0.297477, 0.004030657259573097, 0.004030657259573097, 2, 0.6666666666666666, 0.2167736666666666, 0.6503209999999999, 0.06003183809373913, 0.06003183809373914, 0.021495060013,     0.193455540117,      0.0002419677640180464, 0.0002419677640180465),
    @("FAPs", "L1cam", "Erbb3", "FAPs", 1, 0.3333333333333333, 0.099159, 0.297477, 0.004030657259573097, 0.004030657259573097, 2, 0.6666666666666666, 0.3460976666666666, 1.038293,            0.09584595495126666, 0.09584595495126667, 0.034318698529,     0.308868286761,      0.0003863221941250389, 0.0003863221941250391),
    @("FAPs", "L1cam", "Erbb3", "sCs",  1, 0.3333333333333333, 0.099159, 0.297477, 0.004030657259573097, 0.004030657259573097, 3, 1,                   3.048107,            9.144321,            0.8441222069549942,  0.8441222069549942,  0.302247242013,     2.720225178117,      0.003402367301430011, 0.003402367301430012),
    @("sCs",  "L1cam", "Erbb3", "ECs",  3, 1, 4.781202, 14.343606,                 0.1943483350052483,  0.1943483350052483,  2, 0.6666666666666666, 0.2167736666666666, 0.6503209999999999, 0.06003183809373913, 0.06003183809373914, 1.036438688614,     9.327948197526,      0.01166708778082284, 0.01166708778082284),
    @("sCs",  "L1cam", "Erbb3", "FAPs", 3, 1, 4.781202, 14.343606,                 0.1943483350052483,  0.1943483350052483,  2, 0.6666666666666666, 0.3460976666666666, 1.038293,            0.09584595495126666, 0.09584595495126667, 1.654762856062,     14.892865704558,     0.01862750176176671, 0.01862750176176671),
    @("sCs",  "L1cam", "Erbb3", "sCs",  3, 1, 4.781202, 14.343606,                 0.1943483350052483,  0.1943483350052483,  3, 1,                   3.048107,            9.144321,            0.8441222069549942,  0.8441222069549942,  14.573615284614,    131.162537561526,    0.1640537454626587,  0.1640537454626587)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    for ($c = 5; $c -le 20; $c++) {
        $ws.Cells.Item($row, $c).Value = $r[$c - 1]
    }
    $row++
}
